$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "64.417.54"
$ws.Range("E2").Value = "  -2.40%  "

$ws.Range("D3").Value = "3.343.61"
$ws.Range("E3").Value = "  -4.14%  "

$ws.Range("E4").Value = "  +0.15%  "

$ws.Range("D5").Value = "'554.04"
$ws.Range("E5").Value = "  -4.54%  "

$ws.Range("D6").Value = "'175.28"
$ws.Range("E6").Value = "  -1.41%  "

$ws.Range("E7").Value = "  -2.17%  "

$ws.Range("D8").Value = "3.334.47"
$ws.Range("E8").Value = "  -4.19%  "

$ws.Range("B10").Value = "Dogecoin"
$ws.Range("C10").Value = "https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge"
$ws.Range("D10").Value = "'0.163"
$ws.Range("E10").Value = "  +0.94%  "

$ws.Range("B11").Value = "Cardano"
$ws.Range("C11").Value = "https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada"
$ws.Range("D11").Value = "'0.626"
$ws.Range("E11").Value = "  -1.50%  "

$ws.Range("D12").Value = "'54.44"
$ws.Range("E12").Value = "  -2.63%  "

$ws.Range("D13").Value = "'0.0000273"
$ws.Range("E13").Value = "  -1.90%  "

$ws.Range("E14").Value = "  -2.18%  "

$ws.Range("D15").Value = "3.885.37"
$ws.Range("E15").Value = "  -3.94%  "

$ws.Range("D16").Value = "'18.31"
$ws.Range("E16").Value = "  +0.22%  "

$ws.Range("E17").Value = "  -2.81%  "

$ws.Range("D18").Value = "3.347.11"
$ws.Range("E18").Value = "  -4.03%  "

$ws.Range("D19").Value = "64.374.47"
$ws.Range("E19").Value = "  -2.35%  "

$ws.Range("D20").Value = "'11.74"
$ws.Range("E20").Value = "  -2.60%  "

$ws.Range("D21").Value = "'0.977"
$ws.Range("E21").Value = "  -3.14%  "

$ws.Range("D22").Value = "'432.51"
$ws.Range("E22").Value = "  +5.27%  "

$ws.Range("D23").Value = "'5.09"
$ws.Range("E23").Value = "  +13.53%  "

$ws.Range("E24").Value = "  -4.79%  "

$ws.Range("D25").Value = "'84.24"

$ws.Range("D26").Value = "'13.32"
$ws.Range("E26").Value = "  -0.73%  "

$ws.Range("D27").Value = "'10.75"
$ws.Range("E27").Value = "  -2.60%  "

$ws.Range("D28").Value = "'2.82"
$ws.Range("E28").Value = "  -1.19%  "

$ws.Range("D29").Value = "'8.73"
$ws.Range("E29").Value = "  -4.83%  "

$ws.Range("E30").Value = "  -1.56%  "

$ws.Range("D31").Value = "'6.63"
$ws.Range("E31").Value = "  -0.12%  "

$ws.Range("D32").Value = "'11.48"
$ws.Range("E32").Value = "  -2.05%  "

$ws.Range("D33").Value = "'581.05"
$ws.Range("E33").Value = "  -2.10%  "

$ws.Range("E34").Value = "  -2.87%  "

$ws.Range("D35").Value = "'58.30"
$ws.Range("E35").Value = "  -4.29%  "

$ws.Range("D36").Value = "'0.999"
$ws.Range("E36").Value = "  -0.04%  "

$ws.Range("E37").Value = "  -7.93%  "

$ws.Range("D38").Value = "'3.48"
$ws.Range("E38").Value = "  -3.33%  "

$ws.Range("D39").Value = "'35.65"
$ws.Range("E39").Value = "  -3.10%  "

$ws.Range("D40").Value = "0.0₃0750"
$ws.Range("E40").Value = "  -5.66%  "

$ws.Range("E41").Value = "  -4.50%  "

$ws.Range("D42").Value = "3.112.17"
$ws.Range("E42").Value = "  -3.41%  "

$ws.Range("E43").Value = "  +0.25%  "

$ws.Range("D44").Value = "'2.80"
$ws.Range("E44").Value = "  -5.23%  "

$ws.Range("E45").Value = "  -2.57%  "

$ws.Range("D46").Value = "'0.0408"
$ws.Range("E46").Value = "  -2.35%  "

$ws.Range("D47").Value = "'2.46"
$ws.Range("E47").Value = "  -3.59%  "

$ws.Range("E48").Value = "  -2.06%  "

$ws.Range("E49").Value = "  -3.19%  "

$ws.Range("D50").Value = "'135.54"
$ws.Range("E50").Value = "  -2.04%  "

$ws.Range("D51").Value = "'8.25"
$ws.Range("E51").Value = "  -3.93%  "
